$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '26.699.77'
$ws.Range('E2').Value2 = '  -1.56%  '
$ws.Range('D3').Value2 = '1.794.96'
$ws.Range('E3').Value2 = '  -1.52%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value2 = '0.9999'
$ws.Range('E4').Value2 = '  -0.23%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value2 = '309.09'
$ws.Range('E5').Value2 = '  -0.50%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value2 = '0.9999'
$ws.Range('E6').Value2 = '  -0.13%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value2 = '0.4450'
$ws.Range('E7').Value2 = '  +5.33%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value2 = '0.3671'
$ws.Range('E8').Value2 = '  +0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value2 = '0.07310'
$ws.Range('E9').Value2 = '  +1.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value2 = '0.8577'
$ws.Range('E10').Value2 = '  +0.96%  '
$ws.Range('B11').Value2 = 'Solana'
$ws.Range('C11').Value2 = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value2 = '20.56'
$ws.Range('E11').Value2 = '  -1.71%  '
$ws.Range('B12').Value2 = 'WrappedEther'
$ws.Range('C12').Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value2 = '1.851.33'
$ws.Range('E12').Value2 = '  +0.91%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value2 = '6.624'
$ws.Range('E13').Value2 = '  -0.80%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value2 = '92.59'
$ws.Range('E14').Value2 = '  +3.28%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value2 = '0.07070'
$ws.Range('E15').Value2 = '  -0.39%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value2 = '5.269'
$ws.Range('E16').Value2 = '  -0.41%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value2 = '1.000'
$ws.Range('E17').Value2 = '  -0.30%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value2 = '0.000008666'
$ws.Range('E18').Value2 = '  -1.90%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value2 = '0.9995'
$ws.Range('E19').Value2 = '  -0.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value2 = '14.79'
$ws.Range('E20').Value2 = '  -1.36%  '
$ws.Range('D21').Value2 = '26.725.21'
$ws.Range('E21').Value2 = '  -2.18%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value2 = '5.135'
$ws.Range('E22').Value2 = '  +0.57%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value2 = '10.79'
$ws.Range('E23').Value2 = '  -0.63%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value2 = '1.985'
$ws.Range('E24').Value2 = '  +0.41%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value2 = '151.74'
$ws.Range('E25').Value2 = '  -0.04%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value2 = '18.35'
$ws.Range('E26').Value2 = '  -0.19%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value2 = '2.160'
$ws.Range('E27').Value2 = '  -2.35%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value2 = '5.189'
$ws.Range('E28').Value2 = '  -0.45%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value2 = '117.03'
$ws.Range('E29').Value2 = '  +0.56%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value2 = '0.08773'
$ws.Range('E30').Value2 = '  -0.44%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value2 = '0.7405'
$ws.Range('E31').Value2 = '  -0.63%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value2 = '1.153'
$ws.Range('E32').Value2 = '  -2.91%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value2 = '2.914'
$ws.Range('E33').Value2 = '  -1.74%  '
$ws.Range('E34').Value2 = '  +0.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value2 = '0.9996'
$ws.Range('E35').Value2 = '  -0.16%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value2 = '1.085'
$ws.Range('E36').Value2 = '  -1.55%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value2 = '0.01953'
$ws.Range('E37').Value2 = '  -0.21%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value2 = '0.05174'
$ws.Range('E38').Value2 = '  -1.22%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value2 = '0.5312'
$ws.Range('E39').Value2 = '  +5.56%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value2 = '2.829'
$ws.Range('E40').Value2 = '  -1.62%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value2 = '6.947'
$ws.Range('E41').Value2 = '  -4.72%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value2 = '0.1677'
$ws.Range('E42').Value2 = '  -0.81%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value2 = '0.5076'
$ws.Range('E43').Value2 = '  +7.15%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value2 = '8.403'
$ws.Range('E44').Value2 = '  -2.31%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value2 = '10.48'
$ws.Range('E45').Value2 = '  -0.37%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value2 = '1.935'
$ws.Range('E46').Value2 = '  +3.54%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value2 = '105.18'
$ws.Range('E47').Value2 = '  -1.19%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value2 = '0.9993'
$ws.Range('E48').Value2 = '  -0.14%  '
$ws.Range('E49').Value2 = '  -0.09%  '
$ws.Range('E50').Value2 = '  -1.39%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value2 = '0.9142'
$ws.Range('E51').Value2 = '  +0.76%  '
